# Weekly Fruta/Hortalizas update: a new week of Chirimoya price data
# (date 2023-09-07) is inserted into the "logica_diaria" sheet, just
# after the existing 2022-11-03 block (row 203). This pushes the rest
# of the table (previously rows 203:284) down to rows 206:287, growing
# the used range from A1:T284 to A1:T287.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at 203, shifting rows 203:284 down to 206:287.
$ws.Rows("203:205").Insert()

function Set-ChirimoyaRow {
    param(
        $row,
        $dateSerial,
        $calidad,
        $volumen,
        $precioMin,
        $precioMax,
        $precioProm,
        $unidad,
        $origen,
        $precioKg,
        $kgUnidad
    )
    $ws.Cells.Item($row, 1).Value2 = 8
    $ws.Cells.Item($row, 2).Value = "Terminal La Palmera de La Serena"
    $ws.Cells.Item($row, 3).Value = "Coquimbo"
    $ws.Cells.Item($row, 4).Value2 = $dateSerial
    $ws.Cells.Item($row, 5).Value2 = 4
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value2 = 100107
    $ws.Cells.Item($row, 8).Value = "Otros"
    $ws.Cells.Item($row, 9).Value2 = 100107002
    $ws.Cells.Item($row, 10).Value = "Chirimoya"
    $ws.Cells.Item($row, 11).Value = "Cultivar IV Región"
    $ws.Cells.Item($row, 12).Value = $calidad
    $ws.Cells.Item($row, 13).Value2 = $volumen
    $ws.Cells.Item($row, 14).Value2 = $precioMin
    $ws.Cells.Item($row, 15).Value2 = $precioMax
    $ws.Cells.Item($row, 16).Value2 = $precioProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value2 = $precioKg
    $ws.Cells.Item($row, 20).Value2 = $kgUnidad
}

# New week: 2023-09-07 (serial 45176), Provincia de Limarí, $/bandeja 10 kilos
Set-ChirimoyaRow 203 45176 "Especial" 160 23000 24000 23500 "$/bandeja 10 kilos" "Provincia de Limarí" 2350 10
Set-ChirimoyaRow 204 45176 "Primera"  200 21000 22000 21500 "$/bandeja 10 kilos" "Provincia de Limarí" 2150 10
Set-ChirimoyaRow 205 45176 "Segunda"  200 17000 18000 17500 "$/bandeja 10 kilos" "Provincia de Limarí" 1750 10
